$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (first data row)
$ws.Range("D2").Value = -0.492
$ws.Range("G2").Value = -0.1052631578947368
$ws.Range("H2").Value = -0.1052631578947368
$ws.Range("I2").Value = -0.3578947368421053
$ws.Range("J2").Value = -0.1789473684210526
$ws.Range("K2").Value = -0.887
$ws.Range("L2").Value = -9.336842105263157

$ws.Range("U2").Value = 0.002
$ws.Range("V2").Value = 0.001612903225806452
$ws.Range("W2").Value = 0.05311377245508982
$ws.Range("X2").Value = 0.1016586646398268
$ws.Range("Y2").Value = -0.04854489218473695
$ws.Range("Z2").Value = -0.005707419645539201
$ws.Range("AA2").Value = 0.001021327726043857
$ws.Range("AB2").Value = 0.1011801566806405
$ws.Range("AC2").Value = -0.1001588289545966
$ws.Range("AD2").Value = 0.017
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.017
$ws.Range("AG2").Value = 0.015
$ws.Range("AH2").Value = 0.01352426412092283
$ws.Range("AI2").Value = -0.0009613753322400046
$ws.Range("AJ2").Value = 0.01195219123505976
$ws.Range("AK2").Value = -0.0008481764206955048
$ws.Range("AL2").Value = 0.003
$ws.Range("AM2").Value = -0.011
$ws.Range("AN2").Value = -1.416666666666667
$ws.Range("AO2").Value = -11.33333333333333
$ws.Range("AP2").Value = -1.25
$ws.Range("AQ2").Value = 3.090909090909091

# Row 3 (second data row)
$ws.Range("D3").Value = -0.492
$ws.Range("G3").Value = -0.1052631578947368
$ws.Range("H3").Value = -0.1052631578947368
$ws.Range("I3").Value = -0.3578947368421053
$ws.Range("J3").Value = -0.1789473684210526
$ws.Range("K3").Value = -0.887
$ws.Range("L3").Value = -9.336842105263157

$ws.Range("U3").Value = 0.002
$ws.Range("V3").Value = 0.001612903225806452
$ws.Range("W3").Value = 0.05311377245508982
$ws.Range("X3").Value = 0.1016586646398268
$ws.Range("Y3").Value = -0.04854489218473695
$ws.Range("Z3").Value = -0.005707419645539201
$ws.Range("AA3").Value = 0.001021327726043857
$ws.Range("AB3").Value = 0.1011801566806405
$ws.Range("AC3").Value = -0.1001588289545966
$ws.Range("AD3").Value = 0.017
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.017
$ws.Range("AG3").Value = 0.015
$ws.Range("AH3").Value = 0.01352426412092283
$ws.Range("AI3").Value = -0.0009613753322400046
$ws.Range("AJ3").Value = 0.01195219123505976
$ws.Range("AK3").Value = -0.0008481764206955048
$ws.Range("AL3").Value = 0.003
$ws.Range("AM3").Value = -0.011
$ws.Range("AN3").Value = -1.416666666666667
$ws.Range("AO3").Value = -11.33333333333333
$ws.Range("AP3").Value = -1.25
$ws.Range("AQ3").Value = 3.090909090909091
